# Sprint07.xlsx update:
#  - Closed sprint 7 / started sprint 8: HW3 and "Team Assist" tasks for the
#    last two days of the sprint (Fri/Sat, columns J & K) are now complete,
#    so their remaining-hours burn down to 0 for those days.
#  - Added a note to the Overview comment that the HW3 PDF report is attached
#    and the paper got pushed back again.
#  - Scrolled the Overview sheet back up to the top of the task table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Update the comments/notes cell (A14) with the extra status line.
$ws.Range("A14").Value = "Finish HW3 and get paper drafted.`nHW3 took a little longer than anticipated and so paper got pushed back again."

# HW3 tasks (rows 10-12) and Team Assist (row 13) are finished for the last
# two sprint days (columns J and K) -- clear the carried-forward formula and
# record 0 remaining hours.
$ws.Range("J10:K13").Value = 0

# Scroll the view back up so the selection/top row shown is row 2.
$ws.Application.ActiveWindow.ScrollRow = 2

$wb.Save()
